$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, $tr1.Length).Text = "Testing custom properties"

$tr2 = $s.Shapes.Item(2).TextFrame.TextRange
$tr2.Characters(3, 5).Text = "A. M."
